# Auto-generated Excel COM-interop script
# Applies numeric value updates to the Cerberus_Profits workbook sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the target diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2225.5557
$ws.Range("J17").Value = 2225.5557
$ws.Range("L17").Value = 6676.6671
$ws.Range("N17").Value = -7012.6671

$ws.Range("H62").Value = 1712.1428
$ws.Range("I62").Value = 1712.1428
$ws.Range("K62").Value = 1712.1428
$ws.Range("M62").Value = -1088.1428

$ws.Range("H65").Value = 1712.1428
$ws.Range("I65").Value = 1712.1428
$ws.Range("K65").Value = 8560.714
$ws.Range("M65").Value = -5440.714

$ws.Range("H111").Value = 3128.5
$ws.Range("I111").Value = 2848.5
$ws.Range("J111").Value = 4248.5
$ws.Range("K111").Value = 8545.5
$ws.Range("L111").Value = 12745.5
$ws.Range("M111").Value = -5478.5
$ws.Range("N111").Value = -18879.5

$ws.Range("H113").Value = 6018.913
$ws.Range("I113").Value = 5480.2
$ws.Range("J113").Value = 7029
$ws.Range("K113").Value = 5480.2
$ws.Range("L113").Value = 7029
$ws.Range("M113").Value = -2226.2
$ws.Range("N113").Value = -13537

$ws.Range("H116").Value = 3914.5
$ws.Range("I116").Value = 3926
$ws.Range("J116").Value = 3880
$ws.Range("K116").Value = 3926
$ws.Range("L116").Value = 3880
$ws.Range("M116").Value = -484
$ws.Range("N116").Value = -10764

$ws.Range("H119").Value = 1166.3334
$ws.Range("J119").Value = 1166.3334
$ws.Range("L119").Value = 3499.0002
$ws.Range("N119").Value = -13175.0002

$ws.Range("H129").Value = 1700.5
$ws.Range("I129").Value = 1046.2222
$ws.Range("K129").Value = 3138.6666
$ws.Range("M129").Value = 1861.3334

$ws.Range("H131").Value = 2923.6667

$ws.Range("H137").Value = 1293.8572
$ws.Range("I137").Value = 288
$ws.Range("J137").Value = 1461.5
$ws.Range("K137").Value = 864
$ws.Range("L137").Value = 4384.5
$ws.Range("M137").Value = 1686
$ws.Range("N137").Value = -9484.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 924.95654
$ws.Range("I2").Value = 736.44446
$ws.Range("J2").Value = 1603.6
$ws.Range("K2").Value = 736.44446
$ws.Range("L2").Value = 1603.6
$ws.Range("M2").Value = -623.44446
$ws.Range("N2").Value = -1829.6

$ws.Range("H32").Value = 3185.0645
$ws.Range("I32").Value = 2267.2068
$ws.Range("K32").Value = 2267.2068
$ws.Range("M32").Value = -1980.2068

$ws.Range("H45").Value = 2934.625
$ws.Range("I45").Value = 2020
$ws.Range("K45").Value = 2020
$ws.Range("M45").Value = -1643

$ws.Range("H110").Value = 1338.1052
$ws.Range("I110").Value = 1301.3334
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 1301.3334
$ws.Range("L110").Value = 2000
$ws.Range("M110").Value = 743.6666
$ws.Range("N110").Value = -6090

$ws.Range("H116").Value = 924.95654
$ws.Range("I116").Value = 736.44446
$ws.Range("J116").Value = 1603.6
$ws.Range("K116").Value = 736.44446
$ws.Range("L116").Value = 1603.6
$ws.Range("M116").Value = 1557.55554
$ws.Range("N116").Value = -6191.6

$ws.Range("H132").Value = 1635.4242
$ws.Range("I132").Value = 1547.4193
$ws.Range("K132").Value = 4642.257900000001
$ws.Range("M132").Value = -2112.257900000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 924.95654
$ws.Range("I3").Value = 736.44446
$ws.Range("J3").Value = 1603.6
$ws.Range("K3").Value = 736.44446
$ws.Range("L3").Value = 1603.6
$ws.Range("M3").Value = -622.44446
$ws.Range("N3").Value = -1831.6

$ws.Range("H25").Value = 2077.5
$ws.Range("I25").Value = 2077.5
$ws.Range("K25").Value = 2077.5
$ws.Range("M25").Value = -1842.5

$ws.Range("H105").Value = 3934.625
$ws.Range("I105").Value = 3044.75
$ws.Range("K105").Value = 3044.75
$ws.Range("M105").Value = -1297.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3163.6223
$ws.Range("I31").Value = 2125.3914
$ws.Range("J31").Value = 4249.0454
$ws.Range("K31").Value = 2125.3914
$ws.Range("L31").Value = 4249.0454
$ws.Range("M31").Value = -1830.3914
$ws.Range("N31").Value = -4839.0454

$ws.Range("H34").Value = 3163.6223
$ws.Range("I34").Value = 2125.3914
$ws.Range("J34").Value = 4249.0454
$ws.Range("K34").Value = 2125.3914
$ws.Range("L34").Value = 4249.0454
$ws.Range("M34").Value = -1923.3914
$ws.Range("N34").Value = -4653.0454

$ws.Range("H99").Value = 1891.7693
$ws.Range("I99").Value = 1799.7
$ws.Range("J99").Value = 2198.6667
$ws.Range("K99").Value = 1799.7
$ws.Range("L99").Value = 2198.6667
$ws.Range("M99").Value = -301.7
$ws.Range("N99").Value = -5194.6667

$ws.Range("H126").Value = 1891.7693
$ws.Range("I126").Value = 1799.7
$ws.Range("J126").Value = 2198.6667
$ws.Range("K126").Value = 5399.1
$ws.Range("L126").Value = 6596.000100000001
$ws.Range("M126").Value = -2929.1
$ws.Range("N126").Value = -11536.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 2825201
$ws.Range("I11").Value = 5649991
$ws.Range("K11").Value = 16949973
$ws.Range("M11").Value = -16949833

$ws.Range("H12").Value = 121.86667
$ws.Range("I12").Value = 113.333336
$ws.Range("J12").Value = 127.55556
$ws.Range("K12").Value = 340.000008
$ws.Range("L12").Value = 382.66668
$ws.Range("M12").Value = -167.000008
$ws.Range("N12").Value = -728.66668

$ws.Range("H46").Value = 3824.3
$ws.Range("J46").Value = 4207
$ws.Range("L46").Value = 12621
$ws.Range("N46").Value = -12803

$ws.Range("H57").Value = 7624.75
$ws.Range("I57").Value = 6500
$ws.Range("J57").Value = 10999
$ws.Range("K57").Value = 19500
$ws.Range("L57").Value = 32997
$ws.Range("M57").Value = -18941
$ws.Range("N57").Value = -34115

$ws.Range("H63").Value = 4999
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 4999
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H97").Value = 393
$ws.Range("I97").Value = 393
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1179
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws.Range("H98").Value = 5999
$ws.Range("J98").Value = 5999
$ws.Range("L98").Value = 17997
$ws.Range("N98").Value = -20993

$ws.Range("H107").Value = 1344.375
$ws.Range("I107").Value = 384
$ws.Range("J107").Value = 1920.6
$ws.Range("K107").Value = 1152
$ws.Range("L107").Value = 5761.799999999999
$ws.Range("M107").Value = 768
$ws.Range("N107").Value = -9601.799999999999

$ws.Range("H122").Value = 931.0625
$ws.Range("J122").Value = 1400.8889
$ws.Range("L122").Value = 12608.0001
$ws.Range("N122").Value = -17508.0001

$ws.Range("H131").Value = 1560
$ws.Range("J131").Value = 2366.75
$ws.Range("L131").Value = 7100.25
$ws.Range("N131").Value = -17180.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7253.6665
$ws.Range("I70").Value = 6338.25
$ws.Range("J70").Value = 7711.375
$ws.Range("K70").Value = 6338.25
$ws.Range("L70").Value = 7711.375
$ws.Range("M70").Value = -6068.25
$ws.Range("N70").Value = -8251.375

$ws.Range("H73").Value = 7253.6665
$ws.Range("I73").Value = 6338.25
$ws.Range("J73").Value = 7711.375
$ws.Range("K73").Value = 6338.25
$ws.Range("L73").Value = 7711.375
$ws.Range("M73").Value = -5402.25
$ws.Range("N73").Value = -9583.375

$ws.Range("H126").Value = 8758.9
$ws.Range("I126").Value = 7834.143
$ws.Range("J126").Value = 10916.667
$ws.Range("K126").Value = 23502.429
$ws.Range("L126").Value = 32750.001
$ws.Range("M126").Value = -21032.429
$ws.Range("N126").Value = -37690.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1379
$ws.Range("J46").Value = 1515.1052
$ws.Range("L46").Value = 1515.1052
$ws.Range("N46").Value = -1891.1052

$ws.Range("H61").Value = 3734.7
$ws.Range("I61").Value = 3760.889
$ws.Range("K61").Value = 3760.889
$ws.Range("M61").Value = -3558.889

$ws.Range("H68").Value = 2519.4736
$ws.Range("J68").Value = 2566.3333
$ws.Range("L68").Value = 2566.3333
$ws.Range("N68").Value = -4064.3333

$ws.Range("H71").Value = 2519.4736
$ws.Range("J71").Value = 2566.3333
$ws.Range("L71").Value = 12831.6665
$ws.Range("N71").Value = -20319.6665

$ws.Range("H82").Value = 2355.5
$ws.Range("I82").Value = 2116
$ws.Range("J82").Value = 2535.125
$ws.Range("K82").Value = 2116
$ws.Range("L82").Value = 2535.125
$ws.Range("M82").Value = -1755
$ws.Range("N82").Value = -3257.125

$ws.Range("H85").Value = 2355.5
$ws.Range("I85").Value = 2116
$ws.Range("J85").Value = 2535.125
$ws.Range("K85").Value = 2116
$ws.Range("L85").Value = 2535.125
$ws.Range("M85").Value = -868
$ws.Range("N85").Value = -5031.125

$ws.Range("H113").Value = 3734.7
$ws.Range("I113").Value = 3760.889
$ws.Range("K113").Value = 3760.889
$ws.Range("M113").Value = -1590.889

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7074.5
$ws.Range("I62").Value = 5150
$ws.Range("K62").Value = 5150
$ws.Range("M62").Value = -4526

$ws.Range("H65").Value = 7074.5
$ws.Range("I65").Value = 5150
$ws.Range("K65").Value = 25750
$ws.Range("M65").Value = -22630
